# Adding results of running Budget KP with 10 percent budget on random datasets
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# KPB Training Acc (E) and KPB Test Acc (F) values for data file rows 1..20 (sheet rows 4..23).
# These mirror the KP Training Acc (B) / KP Test Acc (C) values already present.
$values = @(
    @(0.96599999999999997, 0.96299999999999997),
    @(1, 0.997),
    @(0.99399999999999999, 0.99299999999999999),
    @(0.98899999999999999, 0.999),
    @(0.96699999999999997, 0.97799999999999998),
    @(0.98399999999999999, 0.98899999999999999),
    @(0.98799999999999999, 0.98599999999999999),
    @(0.996, 0.996),
    @(0.97399999999999998, 0.97299999999999998),
    @(0.998, 0.999),
    @(0.98199999999999998, 0.97699999999999998),
    @(0.98799999999999999, 0.98199999999999998),
    @(0.98, 0.97299999999999998),
    @(1, 0.99399999999999999),
    @(0.997, 0.998),
    @(0.94499999999999995, 0.94499999999999995),
    @(0.98799999999999999, 0.98699999999999999),
    @(0.997, 0.99299999999999999),
    @(1, 0.997),
    @(0.98899999999999999, 0.997)
)

$row = 4
foreach ($pair in $values) {
    $ws.Cells.Item($row, 5).Value = $pair[0]
    $ws.Cells.Item($row, 6).Value = $pair[1]
    $row++
}

# Update the saved selection to match the author's last-active cell.
$ws.Range("E24").Select()
